$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "26.947.35"
$c.ClearFormats()
$ws.Range("E2").Value = "  +0.41%  "

$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "1.842.47"
$c.ClearFormats()
$ws.Range("E3").Value = "  +0.29%  "

$ws.Range("E4").Value = "  +0.48%  "

$ws.Range("E5").Value = "  +0.37%  "

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "308.75"
$c.ClearFormats()
$ws.Range("E6").Value = "  -0.31%  "

$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.4765"
$c.ClearFormats()
$ws.Range("E7").Value = "  +1.83%  "

$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.3668"
$c.ClearFormats()
$ws.Range("E8").Value = "  +1.23%  "

$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.07208"
$c.ClearFormats()
$ws.Range("E9").Value = "  +0.54%  "

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "0.9292"
$c.ClearFormats()
$ws.Range("E10").Value = "  -0.75%  "

$ws.Range("E11").Value = "  +0.90%  "

$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.07729"
$c.ClearFormats()
$ws.Range("E12").Value = "  +0.64%  "

$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "1.844.52"
$c.ClearFormats()
$ws.Range("E13").Value = "  -0.45%  "

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "5.372"
$c.ClearFormats()
$ws.Range("E14").Value = "  +1.69%  "

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "6.437"
$c.ClearFormats()
$ws.Range("E15").Value = "  +1.02%  "

$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "88.81"
$c.ClearFormats()
$ws.Range("E16").Value = "  +0.90%  "

$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "1.014"
$c.ClearFormats()
$ws.Range("E17").Value = "  +0.49%  "

$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "0.000008656"
$c.ClearFormats()
$ws.Range("E18").Value = "  +0.98%  "

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "1.012"
$c.ClearFormats()
$ws.Range("E19").Value = "  +0.50%  "

$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "27.051.84"
$c.ClearFormats()
$ws.Range("E20").Value = "  +0.78%  "

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "14.52"
$c.ClearFormats()
$ws.Range("E21").Value = "  +1.42%  "

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "5.065"
$c.ClearFormats()
$ws.Range("E22").Value = "  +0.70%  "

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "10.62"
$c.ClearFormats()
$ws.Range("E23").Value = "  +0.03%  "

$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "1.934"
$c.ClearFormats()
$ws.Range("E24").Value = "  +0.88%  "

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "152.67"
$c.ClearFormats()
$ws.Range("E25").Value = "  +0.24%  "

$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "18.17"
$c.ClearFormats()
$ws.Range("E26").Value = "  +0.95%  "

$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "2.004"
$c.ClearFormats()
$ws.Range("E27").Value = "  -0.04%  "

$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "114.22"
$c.ClearFormats()
$ws.Range("E28").Value = "  +0.23%  "

$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "4.964"
$c.ClearFormats()
$ws.Range("E29").Value = "  +1.17%  "

$ws.Range("E30").Value = "  +0.13%  "

$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "3.294"
$c.ClearFormats()
$ws.Range("E31").Value = "  +4.19%  "

$ws.Range("E32").Value = "  -0.85%  "

$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "0.7391"
$c.ClearFormats()
$ws.Range("E33").Value = "  -0.85%  "

$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "4.489"
$c.ClearFormats()
$ws.Range("E34").Value = "  +0.73%  "

$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "2.701"
$c.ClearFormats()
$ws.Range("E35").Value = "  -5.33%  "

$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "1.109"
$c.ClearFormats()
$ws.Range("E36").Value = "  +1.95%  "

$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "0.01957"
$c.ClearFormats()
$ws.Range("E37").Value = "  +1.49%  "

$ws.Range("E38").Value = "  +2.00%  "

$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "2.970"
$c.ClearFormats()
$ws.Range("E39").Value = "  -0.62%  "

$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "0.5237"
$c.ClearFormats()
$ws.Range("E40").Value = "  +2.54%  "

$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "7.009"
$c.ClearFormats()
$ws.Range("E41").Value = "  +1.43%  "

$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "0.1511"
$c.ClearFormats()
$ws.Range("E42").Value = "  -0.25%  "

$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "8.267"
$c.ClearFormats()
$ws.Range("E43").Value = "  +1.53%  "

$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "10.62"
$c.ClearFormats()
$ws.Range("E44").Value = "  +3.41%  "

$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "0.4731"
$c.ClearFormats()
$ws.Range("E45").Value = "  +0.60%  "

$ws.Range("E46").Value = "  +0.38%  "

$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "101.74"
$c.ClearFormats()
$ws.Range("E47").Value = "  +1.96%  "

$ws.Range("E48").Value = "  +1.03%  "

$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "65.76"
$c.ClearFormats()
$ws.Range("E49").Value = "  +2.62%  "

$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "0.06059"
$c.ClearFormats()
$ws.Range("E50").Value = "  +0.28%  "

$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "0.8866"
$c.ClearFormats()
$ws.Range("E51").Value = "  +3.09%  "
